# Updated cryptos list (Price/Volume(1h) refresh + EnergySwap/Monero row swap).
# For numeric-looking price strings we go through a scratch cell holding a
# text formula (="...") and Copy/PasteSpecial(values) into the destination
# so Excel stores them as genuine text (matching the source data, which is
# not true numeric data, e.g. "68.570.39") instead of auto-converting them
# into number cells and losing the original text formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$scratch = $ws.Range("Z1")

$scratch.Formula = '="68.570.39"'
$scratch.Copy()
$ws.Range("D2").PasteSpecial(-4163)
$ws.Range("E2").Value = "  +0.14%  "

$scratch.Formula = '="3.905.45"'
$scratch.Copy()
$ws.Range("D3").PasteSpecial(-4163)
$ws.Range("E3").Value = "  +0.20%  "

$scratch.Formula = '="1.00"'
$scratch.Copy()
$ws.Range("D4").PasteSpecial(-4163)
$ws.Range("E4").Value = "  +0.16%  "

$scratch.Formula = '="602.98"'
$scratch.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  +0.14%  "

$scratch.Formula = '="168.90"'
$scratch.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  +1.59%  "

$scratch.Formula = '="3.903.91"'
$scratch.Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("E8").Value = "  +0.06%  "

$scratch.Formula = '="0.531"'
$scratch.Copy()
$ws.Range("D9").PasteSpecial(-4163)
$ws.Range("E9").Value = "  +0.40%  "

$ws.Range("E10").Value = "  +0.00%  "

$scratch.Formula = '="6.44"'
$scratch.Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = "  +0.49%  "

$scratch.Formula = '="0.460"'
$scratch.Copy()
$ws.Range("D12").PasteSpecial(-4163)
$ws.Range("E12").Value = "  +0.05%  "

$scratch.Formula = '="0.0000253"'
$scratch.Copy()
$ws.Range("D13").PasteSpecial(-4163)
$ws.Range("E13").Value = "  -0.67%  "

$scratch.Formula = '="37.18"'
$scratch.Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  -0.34%  "

$scratch.Formula = '="4.566.36"'
$scratch.Copy()
$ws.Range("D15").PasteSpecial(-4163)
$ws.Range("E15").Value = "  +0.36%  "

$scratch.Formula = '="3.905.58"'
$scratch.Copy()
$ws.Range("D16").PasteSpecial(-4163)
$ws.Range("E16").Value = "  +0.19%  "

$scratch.Formula = '="68.569.89"'
$scratch.Copy()
$ws.Range("D17").PasteSpecial(-4163)
$ws.Range("E17").Value = "  +0.03%  "

$scratch.Formula = '="18.19"'
$scratch.Copy()
$ws.Range("D18").PasteSpecial(-4163)
$ws.Range("E18").Value = "  +5.99%  "

$scratch.Formula = '="7.44"'
$scratch.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = "  -0.40%  "

$ws.Range("E20").Value = "  +0.32%  "

$scratch.Formula = '="10.90"'
$scratch.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  -1.19%  "

$scratch.Formula = '="472.76"'
$scratch.Copy()
$ws.Range("D22").PasteSpecial(-4163)
$ws.Range("E22").Value = "  -2.94%  "

$scratch.Formula = '="0.742"'
$scratch.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = "  +2.50%  "

$ws.Range("E24").Value = "  +0.05%  "

$scratch.Formula = '="83.85"'
$scratch.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = "  -0.77%  "

$ws.Range("E26").Value = "  +1.15%  "

$scratch.Formula = '="12.25"'
$scratch.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  +1.55%  "

$scratch.Formula = '="10.03"'
$scratch.Copy()
$ws.Range("D28").PasteSpecial(-4163)
$ws.Range("E28").Value = "  -0.96%  "

$ws.Range("E29").Value = "  +0.13%  "

$ws.Range("E30").Value = "  +1.20%  "

$scratch.Formula = '="4.057.91"'
$scratch.Copy()
$ws.Range("D31").PasteSpecial(-4163)
$ws.Range("E31").Value = "  +0.22%  "

$ws.Range("E32").Value = "  +1.71%  "

$scratch.Formula = '="31.56"'
$scratch.Copy()
$ws.Range("D33").PasteSpecial(-4163)
$ws.Range("E33").Value = "  -0.74%  "

$ws.Range("E34").Value = "  -2.53%  "

$scratch.Formula = '="9.46"'
$scratch.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = "  +2.01%  "

$scratch.Formula = '="3.883.75"'
$scratch.Copy()
$ws.Range("D36").PasteSpecial(-4163)
$ws.Range("E36").Value = "  +0.86%  "

$ws.Range("E37").Value = "  -1.81%  "

$scratch.Formula = '="3.65"'
$scratch.Copy()
$ws.Range("D38").PasteSpecial(-4163)
$ws.Range("E38").Value = "  +14.85%  "

$scratch.Formula = '="1.03"'
$scratch.Copy()
$ws.Range("D39").PasteSpecial(-4163)
$ws.Range("E39").Value = "  -0.05%  "

$ws.Range("E40").Value = "  +2.48%  "

$ws.Range("E41").Value = "  +0.18%  "

$ws.Range("E42").Value = "  +0.13%  "

$ws.Range("E43").Value = "  -0.20%  "

$scratch.Formula = '="430.30"'
$scratch.Copy()
$ws.Range("D44").PasteSpecial(-4163)
$ws.Range("E44").Value = "  +0.19%  "

$ws.Range("E45").Value = "  +0.98%  "

$ws.Range("E46").Value = "  +13.39%  "

$scratch.Formula = '="8.64"'
$scratch.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = "  +1.41%  "

$scratch.Formula = '="47.29"'
$scratch.Copy()
$ws.Range("D49").PasteSpecial(-4163)
$ws.Range("E49").Value = "  -2.20%  "

# Row 50 becomes EnergySwap, Row 51 becomes Monero (swap with updated values)
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$scratch.Formula = '="26.84"'
$scratch.Copy()
$ws.Range("D50").PasteSpecial(-4163)
$ws.Range("E50").Value = "  +2.96%  "

$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$scratch.Formula = '="143.63"'
$scratch.Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  +1.01%  "

$excel.CutCopyMode = $false
$scratch.Clear()